$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.431.79'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '3.622.41'
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '204.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '568.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.00%  '
$ws.Range("D7").Value = '3.616.32'
$ws.Range("E7").Value = '  +1.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.622'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.59%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.677'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '61.49'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +16.64%  '
$ws.Range("E12").Value = '  +5.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000289'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.35%  '
$ws.Range("D15").Value = '4.187.73'
$ws.Range("E15").Value = '  +1.28%  '
$ws.Range("D16").Value = '3.607.41'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.127'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.20%  '
$ws.Range("D19").Value = '68.201.78'
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '404.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +15.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +16.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '670.89'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.115'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '64.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("E37").Value = '  +2.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.424'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.41%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = '0.0₃0776'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = '3.268.33'
$ws.Range("E41").Value = '  +9.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.87%  '
$ws.Range("E43").Value = '  +3.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.12%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +30.87%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0419'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.132'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.45%  '

Write-Host "Applied all cryptos list updates"